# Add the new "WY Accountability Part 1" mzr_report entry (Scott, Wolfe, Rice, & Wright)
# as the new first item in the mzr_report block, renumbering the existing entries
# that follow it (order column) and shifting them down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renumber the "order" column for the prpa-type rows above the mzr_report block ---
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(4, 2).Value = 3

# --- Insert a new row for the new publication, right before the current first
#     mzr_report row (row 12), pushing the existing mzr_report rows down one ---
$ws.Rows("12:12").Insert()

# Insert() copies the row-11 formatting into the new row 12, including the extra
# P:S tracker cells that don't belong in the mzr_report block - clear those back out.
$ws.Range("P12:S12").Clear()

# --- Populate the new row 12 with the new publication's data ---
$ws.Cells.Item(12, 1).Value = $ws.Cells.Item(13, 1).Value2   # type -> mzr_report (same as rest of block)
$ws.Cells.Item(12, 2).Value = 1                               # order
$ws.Cells.Item(12, 3).Value = "Scott, C., Wolfe, C., Rice, D., & Wright, J."  # authors
$ws.Cells.Item(12, 4).Value = 2023                            # year
$ws.Cells.Item(12, 5).Value = "Wyoming Accountability Part 1: An examination of the current system" # title
$ws.Cells.Item(12, 6).Value = $ws.Cells.Item(13, 6).Value2    # prpa -> Marzano Research (same as rest of block)

# Match the row height Excel computed for the wrapped title text in this row.
$ws.Rows("12:12").RowHeight = 68

# --- Renumber the "order" column for the mzr_report rows that shifted down ---
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(14, 2).Value = 3
$ws.Cells.Item(15, 2).Value = 4
$ws.Cells.Item(16, 2).Value = 5

# --- Update the active selection to match the saved view ---
$ws.Range("C13").Select() | Out-Null
